# Collapse the 3-book listing sheet down to a single "Sacred Tree" listing.
# Row 2 (previously "Two Gentlemen of Verona") is overwritten with the bibliographic
# details that used to live in row 3 ("The Sacred Tree"); the title is replaced with a
# new, shortened form, the price is bumped to 1000, and the old helper column (AN) is
# cleared. Rows 3 and 4 (the Sacred Tree duplicate and the Pearl/Gawain listing) are
# then deleted outright.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Listings")

# Pull row 3's values (photo URL / description / author) across into row 2 before the
# source row disappears. Value2 avoids locale/currency formatting surprises.
$ws.Range("M2").Value = $ws.Range("M3").Value2
$ws.Range("P2").Value = $ws.Range("P3").Value2
$ws.Range("AM2").Value = $ws.Range("AM3").Value2

# New listing title (uses a left single quotation mark, U+2018) + updated start price.
$lq = [char]8216
$ws.Range("E2").Value = "The Sacred Tree, Vol. 2: " + $lq + "The Tale of Genji'"
$ws.Range("K2").Value = 1000

# The redundant "C:Book Title" helper column is no longer populated for this row.
$ws.Range("AN2").ClearContents()

# Drop the old row 3 (duplicate Sacred Tree row) and row 4 (Pearl / Gawain listing).
$ws.Rows("3:4").Delete()

# Restore the window view/selection recorded in the saved workbook.
$ws.Range("AY2").Select()
